$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.386.94"

$ws.Range("D3").Value = "1.923.22"
$ws.Range("E3").Value = "  +3.80%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Formula = "'239.84"
$ws.Range("E5").Value = "  +2.70%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Formula = "'0.4734"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Formula = "'0.2845"
$ws.Range("E8").Value = "  +3.76%  "

$ws.Range("D9").Formula = "'0.06589"
$ws.Range("E9").Value = "  +4.49%  "

$ws.Range("D10").Formula = "'19.11"
$ws.Range("E10").Value = "  +8.09%  "

$ws.Range("D11").Formula = "'104.85"
$ws.Range("E11").Value = "  +24.03%  "

$ws.Range("D12").Value = "1.917.83"
$ws.Range("E12").Value = "  +3.33%  "

$ws.Range("D13").Formula = "'0.07582"
$ws.Range("E13").Value = "  +1.86%  "

$ws.Range("D14").Formula = "'5.121"
$ws.Range("E14").Value = "  +2.86%  "

$ws.Range("D15").Formula = "'0.6510"
$ws.Range("E15").Value = "  +4.04%  "

$ws.Range("D16").Formula = "'298.88"
$ws.Range("E16").Value = "  +21.62%  "

$ws.Range("D17").Value = "30.410.13"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").Formula = "'1.000"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").Formula = "'12.91"
$ws.Range("E19").Value = "  +1.82%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Formula = "'0.000007515"
$ws.Range("E20").Value = "  +2.79%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.162.86"
$ws.Range("E21").Value = "  +2.78%  "

$ws.Range("D22").Formula = "'0.9996"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").Formula = "'5.233"
$ws.Range("E23").Value = "  +5.95%  "

$ws.Range("D24").Formula = "'6.285"
$ws.Range("E24").Value = "  +6.26%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Formula = "'166.53"
$ws.Range("E25").Value = "  +2.35%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Formula = "'9.179"
$ws.Range("E26").Value = "  +0.66%  "

$ws.Range("D27").Formula = "'19.59"
$ws.Range("E27").Value = "  +9.05%  "

$ws.Range("D28").Formula = "'2.031"
$ws.Range("E28").Value = "  +8.44%  "

$ws.Range("D29").Formula = "'0.1117"
$ws.Range("E29").Value = "  +9.37%  "

$ws.Range("D30").Formula = "'1.362"
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").Formula = "'4.098"
$ws.Range("E31").Value = "  +2.32%  "

$ws.Range("D32").Formula = "'3.920"
$ws.Range("E32").Value = "  +2.49%  "

$ws.Range("D33").Formula = "'0.05018"
$ws.Range("E33").Value = "  +3.70%  "

$ws.Range("D34").Formula = "'0.7383"
$ws.Range("E34").Value = "  +5.21%  "

$ws.Range("D35").Formula = "'1.140"
$ws.Range("E35").Value = "  +0.59%  "

$ws.Range("D36").Formula = "'0.9995"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Formula = "'2.715"
$ws.Range("E37").Value = "  +0.52%  "

$ws.Range("D38").Formula = "'0.01949"
$ws.Range("E38").Value = "  +2.84%  "

$ws.Range("D39").Formula = "'2.692"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("E40").Value = "  +2.09%  "

$ws.Range("D41").Formula = "'0.8713"
$ws.Range("E41").Value = "  -0.48%  "

$ws.Range("D42").Formula = "'107.22"

$ws.Range("D43").Formula = "'5.802"
$ws.Range("E43").Value = "  +4.61%  "

$ws.Range("D44").Formula = "'0.9996"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").Formula = "'68.95"
$ws.Range("E45").Value = "  +9.64%  "

$ws.Range("D46").Formula = "'0.4115"
$ws.Range("E46").Value = "  +1.56%  "

$ws.Range("D47").Formula = "'7.243"
$ws.Range("E47").Value = "  +0.86%  "

$ws.Range("D48").Formula = "'9.229"
$ws.Range("E48").Value = "  +8.24%  "

$ws.Range("D49").Formula = "'0.1204"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("D50").Formula = "'34.63"
$ws.Range("E50").Value = "  +3.33%  "

$ws.Range("D51").Formula = "'0.05619"
$ws.Range("E51").Value = "  +1.57%  "
